$wb = $excel.ActiveWorkbook

# "Reinforcement Learning and Opti" is the 4th sheet (sheetId 4 / rId4)
$ws = $wb.Worksheets.Item(4)

$title = "Reinforcement Learning for EV Eco-Driving and Energy Management"
$description = "Task: Optimize the driving policy of an electric vehicle to accurately follow a target speed profile while minimizing energy consumption. The agent must learn to balance speed tracking with efficiency by making intelligent decisions to accelerate, brake, or coast.`nDataset: A custom gymnasium environment using two distinct drive cycles: 1) The standardized Worldwide Harmonised Light Vehicle Test Procedure (WLTP) for training, and 2) Unseen, real-world driving data for validation.`nMethod: Trained a Deep Q-Network (DQN) agent using a custom-engineered reward function. This Reward Shaping technique penalized jerky manoeuvres and incentivized energy-saving coasting to promote a smooth driving style.`nKey Results: The agent successfully learned an efficient, coast-centric policy. The primary achievement was zero-shot generalization: the agent, trained only on the WLTP cycle, applied its efficient driving strategy effectively to the unseen real-world data without any retraining.`nImpact: Demonstrates RL's capability to create robust control policies that can enhance EV range and driving efficiency. This project shows how simulated training can produce intelligent agents ready for complex, real-world automotive challenges."
$link = "https://github.com/Gururaj008/RL_for_EV_Energy_Efficient_Control/"

# Bring formatting for the new row in line with the row above it (row 2)
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)

# Fill in the new project's data
$ws.Range("A3").Value = $title
$ws.Range("B3").Value = $description
$ws.Range("C3").Value = $link

# Let the row height auto-size to the new wrapped content
$ws.Rows.Item(3).AutoFit()

# Turn the URL in C3 into a working hyperlink
$ws.Hyperlinks.Add($ws.Range("C3"), $link)

# Hyperlinks.Add resets the cell style; restore the same "Hyperlink" style used by C2
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# Select / activate this sheet as the active tab, with G2 as the selected cell
$ws.Activate() | Out-Null
$ws.Range("G2").Select() | Out-Null
